$wb = $excel.ActiveWorkbook

# ----- Sheet: Overall -----
$ws = $wb.Worksheets.Item("Overall")

# Header row
$ws.Cells.Item(1, 1).Value = "'Share of 990 filers with government grants at risk"
$ws.Cells.Item(1, 2).Value = "'Number of 990 filers with government grants"
$ws.Cells.Item(1, 3).Value = "'Total government grants (`$)"
$ws.Cells.Item(1, 4).Value = "'Size of operating surplus with government grants"
$ws.Cells.Item(1, 5).Value = "'Size of operating surplus without government grants"

# Data rows
$ws.Cells.Item(2, 1).Value = "'62.94%"
$ws.Cells.Item(2, 2).Value = "'1,592"
$ws.Cells.Item(2, 3).Value = "'`$18,688,408,098"
$ws.Cells.Item(2, 4).Value = "'9.72%"
$ws.Cells.Item(2, 5).Value = "'-9.87%"

# ----- Sheet: County -----
$ws = $wb.Worksheets.Item("County")

# Header row
$ws.Cells.Item(1, 1).Value = "'Geography"
$ws.Cells.Item(1, 2).Value = "'Share of 990 filers with government grants at risk"
$ws.Cells.Item(1, 3).Value = "'Number of 990 filers with government grants"
$ws.Cells.Item(1, 4).Value = "'Total government grants (`$)"
$ws.Cells.Item(1, 5).Value = "'Size of operating surplus with government grants"
$ws.Cells.Item(1, 6).Value = "'Size of operating surplus without government grants"

# Data rows
$ws.Cells.Item(2, 1).Value = "'United States"
$ws.Cells.Item(2, 2).Value = "'67.35%"
$ws.Cells.Item(2, 3).Value = "'103,475"
$ws.Cells.Item(2, 4).Value = "'`$267,700,640,005"
$ws.Cells.Item(2, 5).Value = "'9.05%"
$ws.Cells.Item(2, 6).Value = "'-12.83%"

$ws.Cells.Item(3, 1).Value = "'District of Columbia"
$ws.Cells.Item(3, 2).Value = "'62.94%"
$ws.Cells.Item(3, 3).Value = "'1,592"
$ws.Cells.Item(3, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(3, 5).Value = "'9.72%"
$ws.Cells.Item(3, 6).Value = "'-9.87%"

$ws.Cells.Item(4, 1).Value = "'District of Columbia"
$ws.Cells.Item(4, 2).Value = "'62.94%"
$ws.Cells.Item(4, 3).Value = "'1,592"
$ws.Cells.Item(4, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(4, 5).Value = "'9.72%"
$ws.Cells.Item(4, 6).Value = "'-9.87%"

# ----- Sheet: Congressional District -----
$ws = $wb.Worksheets.Item("Congressional District")

# Header row
$ws.Cells.Item(1, 1).Value = "'Geography"
$ws.Cells.Item(1, 2).Value = "'Share of 990 filers with government grants at risk"
$ws.Cells.Item(1, 3).Value = "'Number of 990 filers with government grants"
$ws.Cells.Item(1, 4).Value = "'Total government grants (`$)"
$ws.Cells.Item(1, 5).Value = "'Size of operating surplus with government grants"
$ws.Cells.Item(1, 6).Value = "'Size of operating surplus without government grants"

# Data rows
$ws.Cells.Item(2, 1).Value = "'United States"
$ws.Cells.Item(2, 2).Value = "'67.35%"
$ws.Cells.Item(2, 3).Value = "'103,475"
$ws.Cells.Item(2, 4).Value = "'`$267,700,640,005"
$ws.Cells.Item(2, 5).Value = "'9.05%"
$ws.Cells.Item(2, 6).Value = "'-12.83%"

$ws.Cells.Item(3, 1).Value = "'District of Columbia"
$ws.Cells.Item(3, 2).Value = "'62.94%"
$ws.Cells.Item(3, 3).Value = "'1,592"
$ws.Cells.Item(3, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(3, 5).Value = "'9.72%"
$ws.Cells.Item(3, 6).Value = "'-9.87%"

$ws.Cells.Item(4, 1).Value = "'Delegate District (at Large)"
$ws.Cells.Item(4, 2).Value = "'62.94%"
$ws.Cells.Item(4, 3).Value = "'1,592"
$ws.Cells.Item(4, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(4, 5).Value = "'9.72%"
$ws.Cells.Item(4, 6).Value = "'-9.87%"

# ----- Sheet: Size -----
$ws = $wb.Worksheets.Item("Size")

# Header row
$ws.Cells.Item(1, 1).Value = "'Size"
$ws.Cells.Item(1, 2).Value = "'Share of 990 filers with government grants at risk"
$ws.Cells.Item(1, 3).Value = "'Number of 990 filers with government grants"
$ws.Cells.Item(1, 4).Value = "'Total government grants (`$)"
$ws.Cells.Item(1, 5).Value = "'Size of operating surplus with government grants"
$ws.Cells.Item(1, 6).Value = "'Size of operating surplus without government grants"

# Data rows
$ws.Cells.Item(2, 1).Value = "'Between `$100K and `$499K"
$ws.Cells.Item(2, 2).Value = "'66.37%"
$ws.Cells.Item(2, 3).Value = "'226"
$ws.Cells.Item(2, 4).Value = "'`$26,963,815"
$ws.Cells.Item(2, 5).Value = "'12.13%"
$ws.Cells.Item(2, 6).Value = "'-19.09%"

$ws.Cells.Item(3, 1).Value = "'Between `$1M and `$4.99M"
$ws.Cells.Item(3, 2).Value = "'61.62%"
$ws.Cells.Item(3, 3).Value = "'594"
$ws.Cells.Item(3, 4).Value = "'`$428,712,597"
$ws.Cells.Item(3, 5).Value = "'10.41%"
$ws.Cells.Item(3, 6).Value = "'-9.20%"

$ws.Cells.Item(4, 1).Value = "'Between `$500K and `$999K"
$ws.Cells.Item(4, 2).Value = "'67.23%"
$ws.Cells.Item(4, 3).Value = "'235"
$ws.Cells.Item(4, 4).Value = "'`$56,142,365"
$ws.Cells.Item(4, 5).Value = "'9.00%"
$ws.Cells.Item(4, 6).Value = "'-12.40%"

$ws.Cells.Item(5, 1).Value = "'Between `$5M and `$9.99M"
$ws.Cells.Item(5, 2).Value = "'51.31%"
$ws.Cells.Item(5, 3).Value = "'191"
$ws.Cells.Item(5, 4).Value = "'`$322,500,981"
$ws.Cells.Item(5, 5).Value = "'12.24%"
$ws.Cells.Item(5, 6).Value = "'-0.87%"

$ws.Cells.Item(6, 1).Value = "'Greater than `$10M"
$ws.Cells.Item(6, 2).Value = "'66.57%"
$ws.Cells.Item(6, 3).Value = "'329"
$ws.Cells.Item(6, 4).Value = "'`$17,853,305,746"
$ws.Cells.Item(6, 5).Value = "'6.84%"
$ws.Cells.Item(6, 6).Value = "'-10.99%"

$ws.Cells.Item(7, 1).Value = "'Less than `$100K"
$ws.Cells.Item(7, 2).Value = "'64.71%"
$ws.Cells.Item(7, 3).Value = "'17"
$ws.Cells.Item(7, 4).Value = "'`$782,594"
$ws.Cells.Item(7, 5).Value = "'40.23%"
$ws.Cells.Item(7, 6).Value = "'-14.85%"

$ws.Cells.Item(8, 1).Value = "'Total"
$ws.Cells.Item(8, 2).Value = "'62.94%"
$ws.Cells.Item(8, 3).Value = "'1,592"
$ws.Cells.Item(8, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(8, 5).Value = "'9.72%"
$ws.Cells.Item(8, 6).Value = "'-9.87%"

# ----- Sheet: Subsector -----
$ws = $wb.Worksheets.Item("Subsector")

# Header row
$ws.Cells.Item(1, 1).Value = "'Subsector"
$ws.Cells.Item(1, 2).Value = "'Share of 990 filers with government grants at risk"
$ws.Cells.Item(1, 3).Value = "'Number of 990 filers with government grants"
$ws.Cells.Item(1, 4).Value = "'Total government grants (`$)"
$ws.Cells.Item(1, 5).Value = "'Size of operating surplus with government grants"
$ws.Cells.Item(1, 6).Value = "'Size of operating surplus without government grants"

# Data rows
$ws.Cells.Item(2, 1).Value = "'Arts, Culture, and Humanities"
$ws.Cells.Item(2, 2).Value = "'73.55%"
$ws.Cells.Item(2, 3).Value = "'155"
$ws.Cells.Item(2, 4).Value = "'`$109,336,570"
$ws.Cells.Item(2, 5).Value = "'7.40%"
$ws.Cells.Item(2, 6).Value = "'-22.94%"

$ws.Cells.Item(3, 1).Value = "'Education (Excluding Universities)"
$ws.Cells.Item(3, 2).Value = "'61.32%"
$ws.Cells.Item(3, 3).Value = "'212"
$ws.Cells.Item(3, 4).Value = "'`$431,916,066"
$ws.Cells.Item(3, 5).Value = "'9.40%"
$ws.Cells.Item(3, 6).Value = "'-8.47%"

$ws.Cells.Item(4, 1).Value = "'Environment and Animals"
$ws.Cells.Item(4, 2).Value = "'63.64%"
$ws.Cells.Item(4, 3).Value = "'77"
$ws.Cells.Item(4, 4).Value = "'`$294,622,837"
$ws.Cells.Item(4, 5).Value = "'9.82%"
$ws.Cells.Item(4, 6).Value = "'-6.54%"

$ws.Cells.Item(5, 1).Value = "'Health (Excluding Hospitals)"
$ws.Cells.Item(5, 2).Value = "'56.76%"
$ws.Cells.Item(5, 3).Value = "'111"
$ws.Cells.Item(5, 4).Value = "'`$267,540,870"
$ws.Cells.Item(5, 5).Value = "'10.52%"
$ws.Cells.Item(5, 6).Value = "'-3.01%"

$ws.Cells.Item(6, 1).Value = "'Hospitals"
$ws.Cells.Item(6, 2).Value = "'100.00%"
$ws.Cells.Item(6, 3).Value = "'3"
$ws.Cells.Item(6, 4).Value = "'`$4,544,365"
$ws.Cells.Item(6, 5).Value = "'19.78%"
$ws.Cells.Item(6, 6).Value = "'-31.83%"

$ws.Cells.Item(7, 1).Value = "'Human Services"
$ws.Cells.Item(7, 2).Value = "'72.70%"
$ws.Cells.Item(7, 3).Value = "'315"
$ws.Cells.Item(7, 4).Value = "'`$989,140,755"
$ws.Cells.Item(7, 5).Value = "'8.74%"
$ws.Cells.Item(7, 6).Value = "'-23.71%"

$ws.Cells.Item(8, 1).Value = "'International, Foreign Affairs"
$ws.Cells.Item(8, 2).Value = "'67.61%"
$ws.Cells.Item(8, 3).Value = "'142"
$ws.Cells.Item(8, 4).Value = "'`$14,237,152,762"
$ws.Cells.Item(8, 5).Value = "'7.68%"
$ws.Cells.Item(8, 6).Value = "'-13.58%"

$ws.Cells.Item(9, 1).Value = "'Mutual/Membership Benefit"
$ws.Cells.Item(9, 2).Value = "'50.00%"
$ws.Cells.Item(9, 3).Value = "'2"
$ws.Cells.Item(9, 4).Value = "'`$1,669,607"
$ws.Cells.Item(9, 5).Value = "'25.74%"
$ws.Cells.Item(9, 6).Value = "'-30.49%"

$ws.Cells.Item(10, 1).Value = "'Public, Societal Benefit"
$ws.Cells.Item(10, 2).Value = "'52.28%"
$ws.Cells.Item(10, 3).Value = "'241"
$ws.Cells.Item(10, 4).Value = "'`$510,468,702"
$ws.Cells.Item(10, 5).Value = "'12.57%"
$ws.Cells.Item(10, 6).Value = "'-1.75%"

$ws.Cells.Item(11, 1).Value = "'Religion Related"
$ws.Cells.Item(11, 2).Value = "'36.36%"
$ws.Cells.Item(11, 3).Value = "'22"
$ws.Cells.Item(11, 4).Value = "'`$42,369,878"
$ws.Cells.Item(11, 5).Value = "'15.09%"
$ws.Cells.Item(11, 6).Value = "'6.87%"

$ws.Cells.Item(12, 1).Value = "'Unclassified"
$ws.Cells.Item(12, 2).Value = "'58.88%"
$ws.Cells.Item(12, 3).Value = "'304"
$ws.Cells.Item(12, 4).Value = "'`$1,106,137,720"
$ws.Cells.Item(12, 5).Value = "'10.60%"
$ws.Cells.Item(12, 6).Value = "'-5.70%"

$ws.Cells.Item(13, 1).Value = "'Universities"
$ws.Cells.Item(13, 2).Value = "'50.00%"
$ws.Cells.Item(13, 3).Value = "'8"
$ws.Cells.Item(13, 4).Value = "'`$693,507,966"
$ws.Cells.Item(13, 5).Value = "'11.95%"
$ws.Cells.Item(13, 6).Value = "'-0.50%"

$ws.Cells.Item(14, 1).Value = "'Total"
$ws.Cells.Item(14, 2).Value = "'62.94%"
$ws.Cells.Item(14, 3).Value = "'1,592"
$ws.Cells.Item(14, 4).Value = "'`$18,688,408,098"
$ws.Cells.Item(14, 5).Value = "'9.72%"
$ws.Cells.Item(14, 6).Value = "'-9.87%"

